$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 0.0004166666666666667
$ws.Range("K2").Value = 2969
$ws.Range("L2").Value = 0.005938
